$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 36.35426266666666
$ws.Range("H2").Value = 109.062788
$ws.Range("I2").Value = 0.4094848412143908
$ws.Range("J2").Value = 0.4094848412143908
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 61.629167
$ws.Range("N2").Value = 184.887501
$ws.Range("O2").Value = 0.8452417044501688
$ws.Range("P2").Value = 0.8452417044501688
$ws.Range("Q2").Value = 2240.482925045865
$ws.Range("R2").Value = 20164.34632541279
$ws.Range("S2").Value = 0.3461136651345584
$ws.Range("T2").Value = 0.3461136651345584

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 36.35426266666666
$ws.Range("H3").Value = 109.062788
$ws.Range("I3").Value = 0.4094848412143908
$ws.Range("J3").Value = 0.4094848412143908
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.357683666666667
$ws.Range("N3").Value = 4.073051
$ws.Range("O3").Value = 0.01862058035795759
$ws.Range("P3").Value = 0.01862058035795759
$ws.Range("Q3").Value = 49.3575886362431
$ws.Range("R3").Value = 444.2182977261879
$ws.Range("S3").Value = 0.007624845391198068
$ws.Range("T3").Value = 0.007624845391198068

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 36.35426266666666
$ws.Range("H4").Value = 109.062788
$ws.Range("I4").Value = 0.4094848412143908
$ws.Range("J4").Value = 0.4094848412143908
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 9.926218666666665
$ws.Range("N4").Value = 29.778656
$ws.Range("O4").Value = 0.1361377151918736
$ws.Range("P4").Value = 0.1361377151918736
$ws.Range("Q4").Value = 360.8603606947697
$ws.Range("R4").Value = 3247.743246252928
$ws.Range("S4").Value = 0.05574633068863431
$ws.Range("T4").Value = 0.05574633068863432

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 45.11545066666667
$ws.Range("H5").Value = 135.346352
$ws.Range("I5").Value = 0.5081685556916724
$ws.Range("J5").Value = 0.5081685556916724
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 61.629167
$ws.Range("N5").Value = 184.887501
$ws.Range("O5").Value = 0.8452417044501688
$ws.Range("P5").Value = 0.8452417044501688
$ws.Range("Q5").Value = 2780.427643416262
$ws.Range("R5").Value = 25023.84879074636
$ws.Range("S5").Value = 0.4295252561608097
$ws.Range("T5").Value = 0.4295252561608097

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 45.11545066666667
$ws.Range("H6").Value = 135.346352
$ws.Range("I6").Value = 0.5081685556916724
$ws.Range("J6").Value = 0.5081685556916724
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.357683666666667
$ws.Range("N6").Value = 4.073051
$ws.Range("O6").Value = 0.01862058035795759
$ws.Range("P6").Value = 0.01862058035795759
$ws.Range("Q6").Value = 61.25251048443911
$ws.Range("R6").Value = 551.272594359952
$ws.Range("S6").Value = 0.009462393426644032
$ws.Range("T6").Value = 0.009462393426644032

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 45.11545066666667
$ws.Range("H7").Value = 135.346352
$ws.Range("I7").Value = 0.5081685556916724
$ws.Range("J7").Value = 0.5081685556916724
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 9.926218666666665
$ws.Range("N7").Value = 29.778656
$ws.Range("O7").Value = 0.1361377151918736
$ws.Range("P7").Value = 0.1361377151918736
$ws.Range("Q7").Value = 447.8258285625458
$ws.Range("R7").Value = 4030.432457062912
$ws.Range("S7").Value = 0.06918090610421863
$ws.Range("T7").Value = 0.06918090610421865

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 7.310771333333334
$ws.Range("H8").Value = 21.932314
$ws.Range("I8").Value = 0.0823466030939367
$ws.Range("J8").Value = 0.0823466030939367
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 61.629167
$ws.Range("N8").Value = 184.887501
$ws.Range("O8").Value = 0.8452417044501688
$ws.Range("P8").Value = 0.8452417044501688
$ws.Range("Q8").Value = 450.5567474008128
$ws.Range("R8").Value = 4055.010726607315
$ws.Range("S8").Value = 0.06960278315480059
$ws.Range("T8").Value = 0.06960278315480059

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 7.310771333333334
$ws.Range("H9").Value = 21.932314
$ws.Range("I9").Value = 0.0823466030939367
$ws.Range("J9").Value = 0.0823466030939367
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.357683666666667
$ws.Range("N9").Value = 4.073051
$ws.Range("O9").Value = 0.01862058035795759
$ws.Range("P9").Value = 0.01862058035795759
$ws.Range("Q9").Value = 9.925714830001555
$ws.Range("R9").Value = 89.33143347001399
$ws.Range("S9").Value = 0.001533341540115487
$ws.Range("T9").Value = 0.001533341540115487

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 7.310771333333334
$ws.Range("H10").Value = 21.932314
$ws.Range("I10").Value = 0.0823466030939367
$ws.Range("J10").Value = 0.0823466030939367
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 9.926218666666665
$ws.Range("N10").Value = 29.778656
$ws.Range("O10").Value = 0.1361377151918736
$ws.Range("P10").Value = 0.1361377151918736
$ws.Range("Q10").Value = 72.56831487666489
$ws.Range("R10").Value = 653.114833889984
$ws.Range("S10").Value = 0.01121047839902061
$ws.Range("T10").Value = 0.01121047839902061
